# The DeviceList sheet had a device column ("APPLE_iPhone13_iOS_15.2.0_fb8f5")
# removed entirely (column F). Deleting the whole column shifts G/H/I/J left
# into F/G/H/I and drops the old column J, which is exactly what the target
# diff shows (dimension A1:J10 -> A1:I10, every row's span 1:10 -> 1:9, the
# shared-strings table losing the two now-unreferenced strings, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

$ws.Columns("F").Delete()

# The two conditional-formatting rule groups on row 2 need their ranges (and,
# for the right-hand group, their formulas) nudged to follow the shifted
# columns: B2:H2 -> B2:G2, and I2:J2 -> H2:I2 with the formula's "I2"
# reference becoming "H2".
$leftRules = $ws.Range("B2:H2").FormatConditions
$leftRules.Item(1).ModifyAppliesToRange($ws.Range("B2:G2"))
$leftRules.Item(2).ModifyAppliesToRange($ws.Range("B2:G2"))

$rightRules = $ws.Range("I2:J2").FormatConditions
$rightRules.Item(1).ModifyAppliesToRange($ws.Range("H2:I2"))
$rightRules.Item(1).Formula1 = '=LEFT((H2),LEN("Working"))=("Working")'
$rightRules.Item(2).ModifyAppliesToRange($ws.Range("H2:I2"))
$rightRules.Item(2).Formula1 = '=LEFT((H2),LEN("Not Working"))=("Not Working")'

# The saved view had scrolled so D1 was the top-left visible cell with J15
# selected; the edited workbook instead shows the default top-left corner
# with D12 selected.
$ws.Range("D12").Select()
